$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-47: recompute columns B (time) and C (Energy)
# using the new priority-queue based simulation results.
$updated = New-Object "object[,]" 46,2
$updated[0,0] = 1.070094600866085
$updated[0,1] = 1.949439117313653
$updated[1,0] = 1.436643150133418
$updated[1,1] = 4.052885245777932
$updated[2,0] = 3.787215586537338
$updated[2,1] = 6.519892421076485
$updated[3,0] = 4.205115307607776
$updated[3,1] = 8.521995367874593
$updated[4,0] = 5.884771137193976
$updated[4,1] = 10.76484375443066
$updated[5,0] = 6.737731502173306
$updated[5,1] = 12.82408792781244
$updated[6,0] = 8.269121302527745
$updated[6,1] = 15.09184207979542
$updated[7,0] = 9.821478568297186
$updated[7,1] = 17.49637731023616
$updated[8,0] = 12.63574511001238
$updated[8,1] = 19.54618568192062
$updated[9,0] = 13.01577696911315
$updated[9,1] = 21.62524704637249
$updated[10,0] = 15.81653229041376
$updated[10,1] = 24.50956169666455
$updated[11,0] = 17.20130088917745
$updated[11,1] = 26.66492695702541
$updated[12,0] = 18.15490055164808
$updated[12,1] = 28.61459237479785
$updated[13,0] = 21.07259023384928
$updated[13,1] = 30.65123655041259
$updated[14,0] = 22.26389876482664
$updated[14,1] = 32.8685040839992
$updated[15,0] = 22.42468537696624
$updated[15,1] = 34.76960797993618
$updated[16,0] = 25.6348929321804
$updated[16,1] = 36.99349727634269
$updated[17,0] = 26.71529546957016
$updated[17,1] = 39.29025016631252
$updated[18,0] = 28.56097429342883
$updated[18,1] = 41.33036300388491
$updated[19,0] = 36.33386712404639
$updated[19,1] = 43.63655750256484
$updated[20,0] = 38.83332052310951
$updated[20,1] = 45.60072540448348
$updated[21,0] = 39.60382002608137
$updated[21,1] = 47.79310263464475
$updated[22,0] = 40.84118351788496
$updated[22,1] = 50.05359660531692
$updated[23,0] = 42.05722374026853
$updated[23,1] = 52.0510276745283
$updated[24,0] = 42.18905378045974
$updated[24,1] = 54.29446044731539
$updated[25,0] = 45.39191774621054
$updated[25,1] = 56.2491001439523
$updated[26,0] = 50.3231234509356
$updated[26,1] = 58.59907624701708
$updated[27,0] = 50.77322372827352
$updated[27,1] = 60.67271899913607
$updated[28,0] = 54.75662843798163
$updated[28,1] = 62.55422647647692
$updated[29,0] = 56.1324871742077
$updated[29,1] = 64.65001942580884
$updated[30,0] = 58.07486926666923
$updated[30,1] = 66.63248690719456
$updated[31,0] = 59.86494961767859
$updated[31,1] = 68.51520602420244
$updated[32,0] = 63.62700101341061
$updated[32,1] = 70.46447443858833
$updated[33,0] = 65.26627530769346
$updated[33,1] = 72.60200273261145
$updated[34,0] = 68.66948468911215
$updated[34,1] = 74.9873895309416
$updated[35,0] = 70.48856712933821
$updated[35,1] = 77.13812851229274
$updated[36,0] = 72.60794727941862
$updated[36,1] = 79.41023042173869
$updated[37,0] = 74.2425974865932
$updated[37,1] = 81.43170463317584
$updated[38,0] = 76.06365240755424
$updated[38,1] = 83.50621878210691
$updated[39,0] = 76.14957874210232
$updated[39,1] = 85.4884105376415
$updated[40,0] = 77.87213565580082
$updated[40,1] = 87.79584633752843
$updated[41,0] = 82.83079939216451
$updated[41,1] = 89.80655793631423
$updated[42,0] = 84.98897450903945
$updated[42,1] = 92.06610034862747
$updated[43,0] = 87.15851009208343
$updated[43,1] = 94.28328766620751
$updated[44,0] = 87.29352796289311
$updated[44,1] = 96.47890792568111
$updated[45,0] = 89.15531203860846
$updated[45,1] = 98.74761496421334

$ws.Range("B2:C47").Value = $updated

# Append three additional simulated rows (46, 47, 48) at the end of the table
$newRows = @(
    @(46, 90.86454508998294, 100.997633300931),
    @(47, 92.63513400026071, 102.9924907112505),
    @(48, 94.92473796977741, 104.9403476662315),
)

$lastRow = 47
foreach ($rowVals in $newRows) {
    $lastRow++
    # Copy the style of the row above (bold/bordered A column) onto the new A cell
    $ws.Range("A" + ($lastRow - 1)).Copy($ws.Range("A" + $lastRow))
    $ws.Range("A" + $lastRow).Value = $rowVals[0]
    $ws.Range("B" + $lastRow).Value = $rowVals[1]
    $ws.Range("C" + $lastRow).Value = $rowVals[2]
}
